$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row labels: "_old" -> "_FV2210", "_new" -> "_FV2304"
#    (columns A..J are the FV2210 side, K is "diff", L..U are the FV2304 side)
# ---------------------------------------------------------------------------
$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") spanning A1:U82
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $dataRange, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, pane = A2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
